$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# 1) Shared string fix: "use 400" -> "CorTec400". Use Find/Replace so every
#    cell referencing the shared string gets updated in place instead of
#    forking a brand-new string for a single cell.
# -------------------------------------------------------------------------
$ws.Cells.Replace("use 400", "CorTec400") | Out-Null

# -------------------------------------------------------------------------
# 2) Row 3: align E3/G3 formatting (fill) with the rest of the table
#    (rows 4-12 already use this fill for columns E/G).
# -------------------------------------------------------------------------
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("G4").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null

# -------------------------------------------------------------------------
# 3) Move the "CorTec400" flag from column L into column H for rows 3-12
#    (row 8 never had one), matching the fill used by E/G on each row.
# -------------------------------------------------------------------------
foreach ($r in 3,4,5,6,7,9,10,11,12) {
  $ws.Cells.Item($r, 5).Copy() | Out-Null
  $ws.Cells.Item($r, 8).PasteSpecial(-4122) | Out-Null
  $ws.Cells.Item($r, 8).Value2 = $ws.Cells.Item($r, 12).Value2
  $ws.Cells.Item($r, 12).ClearContents()
}

# -------------------------------------------------------------------------
# 4) Row 8 gets a brand-new (empty) H8 cell with red font color -- this is
#    the new font (fontId 2 / red) introduced in styles.xml.
# -------------------------------------------------------------------------
$ws.Range("H8").Font.Color = 255

# -------------------------------------------------------------------------
# 5) Column widths: E/F gain explicit widths, H gets one too (mirrors the
#    bestFit columns elsewhere in the sheet).
# -------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 7.307291666666667
$ws.Columns("F").ColumnWidth = 9.166666666666666
$ws.Columns("H").ColumnWidth = 9.166666666666666

# -------------------------------------------------------------------------
# 6) Selection moves from L23 to J24.
# -------------------------------------------------------------------------
$ws.Range("J24").Select() | Out-Null
